$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @"
2|60.032.23|  +4.89%  
3|2.590.60|  +6.62%  
4||  -0.11%  
5|504.71|  +3.26%  
6|155.83|  +1.39%  
7|0.995|  -0.24%  
8|0.588|  -3.63%  
9|2.627.35|  +7.92%  
10|6.45|  +6.13%  
11||  +4.01%  
12|0.341|  +2.51%  
13||  +0.72%  
14|3.069.61|  +8.00%  
15|60.188.68|
16|21.59|  +4.83%  
17|0.0000139|  +4.36%  
18|2.624.58|  +7.98%  
19|4.76|  +2.19%  
20|342.62|  +5.63%  
21|10.41|  +4.10%  
22|6.14|  +3.38%  
23|0.996|  -0.37%  
24|59.84|  +3.42%  
25|0.421|  +4.43%  
26|2.726.62|  +8.02%  
27|0.164|  +2.40%  
28||  -0.51%  
29|0.0₃0853|  +8.52%  
30|7.52|  +3.14%  
31|0.998|  -0.17%  
32|19.42|  +4.26%  
33|155.85|  +3.02%  
34||  +2.94%  
35|5.71|  +7.93%  
36|4.01|  +6.43%  
37||  +4.11%  
38|306.73|  +7.41%  
39|1.48|  +7.83%  
40|0.843|  +3.17%  
41|3.75|  +6.73%  
42|0.834|  +27.91%  
43|35.46|  +4.39%  
44||  +5.45%  
45|0.0570|  +7.18%  
46||  -1.31%  
47|0.993|  -0.33%  
48|19.73|  +11.84%  
49|4.86|  +7.31%  
50|2.044.19|  +7.25%  
51|0.0234|  +2.86%  
"@

$lines = $rowData -split "`n"
foreach ($line in $lines) {
    $line = $line.TrimEnd("`r")
    if ($line -eq "") { continue }
    $parts = $line -split '\|', 3
    $rowNum = [int]$parts[0]
    $dVal = $parts[1]
    $eVal = $parts[2]

    if ($dVal -ne "") {
        if ($dVal -match '^-?\d+(\.\d+)?$') {
            $ws.Cells.Item($rowNum, 4).Value = "'" + $dVal
        } else {
            $ws.Cells.Item($rowNum, 4).Value = $dVal
        }
    }
    if ($eVal -ne "") {
        $ws.Cells.Item($rowNum, 5).Value = $eVal
    }
}
